$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.964.81'
$ws.Range("E2").Value = '  -3.27%  '
# Row 3
$ws.Range("D3").Value = '3.321.71'
$ws.Range("E3").Value = '  -4.79%  '
# Row 4
$ws.Range("E4").Value = '  +0.02%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '182.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -7.40%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '531.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.48%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.605'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.74%  '
# Row 8
$ws.Range("D8").Value = '3.316.56'
$ws.Range("E8").Value = '  -4.52%  '
# Row 9
$ws.Range("E9").Value = '  +0.04%  '
# Row 10
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.620'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.41%  '
# Row 11
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.16'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.56%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.134'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.34%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000262'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.77%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.59%  '
# Row 15
$ws.Range("D15").Value = '3.867.32'
$ws.Range("E15").Value = '  -4.69%  '
# Row 16
$ws.Range("D16").Value = '3.333.67'
$ws.Range("E16").Value = '  -4.50%  '
# Row 17
$ws.Range("E17").Value = '  -4.50%  '
# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.84%  '
# Row 19
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '65.034.14'
$ws.Range("E19").Value = '  -2.84%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.17%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.966'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.27%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '376.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.04%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.50%  '
# Row 24
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.61'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.26%  '
# Row 25
$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.13%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.45%  '
# Row 27
$ws.Range("E27").Value = '  -0.76%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.59%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.27%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.19%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '29.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.04%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '655.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.17%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.08%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.37'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.15%  '
# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.106'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.88%  '
# Row 36
$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '59.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.72%  '
# Row 37
$ws.Range("E37").Value = '  -0.10%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.394'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.68%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.44%  '
# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.24%  '
# Row 41
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '0.0₃0712'
$ws.Range("E41").Value = '  +6.23%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.127'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.12%  '
# Row 43
$ws.Range("D43").Value = '2.913.07'
$ws.Range("E43").Value = '  -4.48%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.02%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.72'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.37%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0403'
$ws.Range("D46").Style = "Normal"
# Row 47
$ws.Range("E47").Value = '  -1.82%  '
# Row 48
$ws.Range("E48").Value = '  +12.06%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.05%  '
# Row 50
$ws.Range("E50").Value = '  +0.65%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.06%  '
